# Auto-generated edit script applying cryptos.xlsx diff (price/volume refresh + row reorder for rows 45-48)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.880.16"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "3.251.48"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'546.80"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").Value = "'149.08"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.525"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'7.47"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "'0.434"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "3.800.35"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "'26.53"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "60.839.72"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "3.252.31"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "'6.37"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").Value = "'13.50"
$ws.Range("E19").Value = "  +3.80%  "
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").Value = "'377.47"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "'0.531"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'69.95"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").Value = "'8.70"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  +6.23%  "
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").Value = "'22.63"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("E33").Value = "  +6.80%  "
$ws.Range("D34").Value = "'6.68"
$ws.Range("E34").Value = "  +5.16%  "
$ws.Range("D35").Value = "'159.99"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("D37").Value = "'26.63"
$ws.Range("E37").Value = "  +3.32%  "
$ws.Range("D38").Value = "2.801.06"
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("D39").Value = "'0.0721"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'1.74"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").Value = "'0.0316"
$ws.Range("E41").Value = "  +7.84%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "'40.18"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.106"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.01"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'21.70"
$ws.Range("E47").Value = "  +6.67%  "
$ws.Range("B48").Value = "RenzoRestakedETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D48").Value = "3.286.57"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").Value = "'0.808"
$ws.Range("E50").Value = "  +6.11%  "
$ws.Range("D51").Value = "'281.26"
$ws.Range("E51").Value = "  +9.42%  "

Write-Host "Applied cryptos.xlsx update"
